# fix bug ru council
# Updates recomputed 'mean' values in the solidarity_support_mean sheet
# after correcting the underlying Russia / Security-Council data bug.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.744386881121362
$ws.Range("L2").Value = 0.650516678087879
$ws.Range("B3").Value = 0.539519237088218
$ws.Range("L3").Value = 0.62733909707771
$ws.Range("B4").Value = 0.458531133781306
$ws.Range("L4").Value = 0.701142356844379
$ws.Range("B5").Value = 0.336663915685529
$ws.Range("C5").Value = 0.340947849011829
$ws.Range("D5").Value = 0.267827509596419
$ws.Range("E5").Value = 0.266127934188177
$ws.Range("F5").Value = 0.513159966971843
$ws.Range("G5").Value = 0.0950055489582969
$ws.Range("H5").Value = 0.60555401572668
$ws.Range("I5").Value = 0.323488971085099
$ws.Range("J5").Value = 0.25889276224129
$ws.Range("K5").Value = 0.0939565076587604
$ws.Range("L5").Value = 0.636887453022506
$ws.Range("M5").Value = 0.808599687276629
$ws.Range("N5").Value = 0.269914170603018
$ws.Range("B6").Value = 0.333653647468452
$ws.Range("C6").Value = 0.366925945043564
$ws.Range("D6").Value = 0.218677332127399
$ws.Range("E6").Value = 0.135982701351322
$ws.Range("F6").Value = 0.664560455495347
$ws.Range("G6").Value = 0.437994926973674
$ws.Range("H6").Value = 0.496762185332231
$ws.Range("I6").Value = 0.441736751387392
$ws.Range("J6").Value = 0.29213655766817
$ws.Range("K6").Value = 0.202429031291367
$ws.Range("L6").Value = 0.373834833505223
$ws.Range("M6").Value = 0.932399297563257
$ws.Range("N6").Value = 0.292154611711296
$ws.Range("B7").Value = 0.328352351803956
$ws.Range("C7").Value = 0.523687810073323
$ws.Range("D7").Value = 0.488432545266734
$ws.Range("E7").Value = 0.463528415737963
$ws.Range("F7").Value = 0.896310327588646
$ws.Range("G7").Value = 0.223242241830577
$ws.Range("H7").Value = 0.461373353988713
$ws.Range("I7").Value = 0.519323835787545
$ws.Range("J7").Value = 0.25978725026746
$ws.Range("K7").Value = 0.331128776431993
$ws.Range("L7").Value = -0.080025210199837
$ws.Range("M7").Value = 0.497163901847429
$ws.Range("N7").Value = 0.266414226280647
$ws.Range("B8").Value = 0.321059513541112
$ws.Range("C8").Value = 0.492563434872145
$ws.Range("D8").Value = 0.382508576398462
$ws.Range("E8").Value = 0.466504525543576
$ws.Range("F8").Value = 0.675638317174231
$ws.Range("G8").Value = 0.340965649884362
$ws.Range("H8").Value = 0.539782540352986
$ws.Range("I8").Value = 0.498286884195477
$ws.Range("J8").Value = 0.40908858490456
$ws.Range("K8").Value = 0.229930305720341
$ws.Range("L8").Value = 0.0107616800013773
$ws.Range("M8").Value = 0.707615721886186
$ws.Range("N8").Value = 0.262541406776756
$ws.Range("B9").Value = 0.31841123807984
$ws.Range("L9").Value = 0.369383378979967
$ws.Range("B10").Value = 0.317414800687998
$ws.Range("L10").Value = 0.673446690200645
$ws.Range("B11").Value = 0.0133623665043539
$ws.Range("C11").Value = 0.0839217421719148
$ws.Range("L11").Value = 0.00385878470855154
